# Rescale the x_m (E) and y_m (F) sensor coordinates to match the
# corrected model calibration. Every value in E2:F61 is multiplied by
# the same constant factor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scale = 0.9968730456535334

$lastRow = 61
for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $eOld = $eCell.Value()
    $fOld = $fCell.Value()
    $eCell.Value = $eOld * $scale
    $fCell.Value = $fOld * $scale
}

Write-Output "Rescaled E2:F61 by factor $scale"
